$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230344891548157
$ws.Range("B1").Value = 2.474378824234009
$ws.Range("C1").Value = 4.022601127624512
$ws.Range("D1").Value = 2.774567127227783
$ws.Range("E1").Value = 1.088649272918701
